$d = $word.ActiveDocument

function Replace-Paragraph($paragraph, [string]$innerXml) {
    $full = $d.Range($paragraph.Range.Start, $paragraph.Range.End)
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData>' +
        '</pkg:part>' +
        '</pkg:package>'
    $full.InsertXML($xml) | Out-Null
}

function Replace-Range($range, [string]$innerXml) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData>' +
        '</pkg:part>' +
        '</pkg:package>'
    $range.InsertXML($xml) | Out-Null
}

# ----------------------------------------------------------------------
# Paragraph 3: intro paragraph rewrite
# ----------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$inner3 = '<w:p>' +
    '<w:r><w:t xml:space="preserve">Voi realiza o aplicatie web unde un utilizator va putea </w:t></w:r>' +
    '<w:r><w:t>vizualiza lista celor mai bune</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> filme, pentru ca utilizatorul sa poata </w:t></w:r>' +
    '<w:r><w:t>accesa lista de filme</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">, el va trebui sa se logheze, daca nu are cont va trebui sa-si creeze unul, deasemenea va putea </w:t></w:r>' +
    '<w:r><w:t>sa adauge un film nou in lista in cazul in care crede ca filmul sau preferat nu este in lista.</w:t></w:r>' +
    '</w:p>'
Replace-Paragraph $p3 $inner3

# ----------------------------------------------------------------------
# Paragraph 4: "Pentru a realiza functionalitatile..." - spring -> spring boot framework
# ----------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$inner4 = '<w:p>' +
    '<w:r><w:t>Pentru a realiza functionalita</w:t></w:r>' +
    '<w:r><w:t>t</w:t></w:r>' +
    '<w:r><w:t>ile descrise mai sus voi crea un proiec</w:t></w:r>' +
    '<w:r><w:t>t</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> nou in InteliiJ IDEA folosind spring</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> boot</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> framework</w:t></w:r>' +
    '<w:r><w:t>:</w:t></w:r>' +
    '</w:p>'
Replace-Paragraph $p4 $inner4

# ----------------------------------------------------------------------
# Paragraph 5: list item about baza de date - ", date privind metoda de plata)" -> ", lista de filme)"
# ----------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$pPr5 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$inner5 = '<w:p>' + $pPr5 +
    '<w:r><w:t>voi realiza conexiunea la baza de dat</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">e cu ajutorul layer-ului model, </w:t></w:r>' +
    '<w:r><w:t>in baza de date voi stoca date despre utilizator(date logare</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
    '<w:r><w:t>lista de filme</w:t></w:r>' +
    '<w:r><w:t>)</w:t></w:r>' +
    '</w:p>'
Replace-Paragraph $p5 $inner5

# ----------------------------------------------------------------------
# Paragraph 7: "Pentru ca utilizatorul sa poata viziona..." rewrite, with
# the _GoBack bookmark now living before the trailing period.
# ----------------------------------------------------------------------
$p7 = $d.Paragraphs.Item(7)
$inner7 = '<w:p>' +
    '<w:r><w:t xml:space="preserve">Pentru ca utilizatorul sa poata </w:t></w:r>' +
    '<w:r><w:t>accesa lista de filme voi</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> folosi metoda GET, deoarece utilizatorul</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> trebuie sa primeasca niste date</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '</w:p>'
Replace-Paragraph $p7 $inner7

# ----------------------------------------------------------------------
# Paragraph 8: "Pentru crearea unul cont nou sau logare..." rewrite
# ----------------------------------------------------------------------
$p8 = $d.Paragraphs.Item(8)
$inner8 = '<w:p>' +
    '<w:r><w:t>Pentru crearea unul cont nou voi folosi metoda POST unde in body voi trimite catre server datele introduse de utilizator</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> si pentru logare voi folosi dependita spring boot starter security</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '</w:p>'
Replace-Paragraph $p8 $inner8

# ----------------------------------------------------------------------
# Paragraphs 9-12: "Pentru cautarea...", "Pentru adaugarea de
# comentariu...", "Pentru adaugarea metodei de plata..." (old bookmark)
# and the picture paragraph all collapse into a single new paragraph.
# ----------------------------------------------------------------------
$p9 = $d.Paragraphs.Item(9)
$p12 = $d.Paragraphs.Item(12)
$range9to12 = $d.Range($p9.Range.Start, $p12.Range.End)
$inner9 = '<w:p>' +
    '<w:r><w:t xml:space="preserve">Pentru </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">adaugarea unui film nou pe lista </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">voi folosi la fel metoda </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">POST </w:t></w:r>' +
    '<w:r><w:t>pentru a transmite datele catre server fiind ulterior salvate intr-o baza de date.</w:t></w:r>' +
    '</w:p>'
Replace-Range $range9to12 $inner9
